# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
# Update Front axle (row 5/6) and Rear axle (row 9/10) hardpoint
# sBottom/xPreload values (columns G/H, and F for rows 9/10), switching
# them from the 3-decimal custom format to a 2-decimal format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Front, sTop "K") ---
$ws.Range("G5").Value = 0.62
$ws.Range("H5").Value = 0.65
$ws.Range("G5:H5").NumberFormat = "0.00"

# --- Row 6 (Front, sBottom "K") ---
$ws.Range("G6").Value = 0.85
$ws.Range("H6").Value = 0.19
$ws.Range("G6:H6").NumberFormat = "0.00"

# --- Row 9 (Rear, sTop "K") ---
$ws.Range("G9").Value = 0.62
$ws.Range("H9").Value = 0.65
$ws.Range("F9:H9").NumberFormat = "0.00"

# --- Row 10 (Rear, sBottom "K") ---
$ws.Range("G10").Value = 0.85
$ws.Range("H10").Value = 0.19
$ws.Range("F10:H10").NumberFormat = "0.00"

# --- Sheet tab color (theme accent5, darker tint -> RGB 2E75B6) ---
$ws.Tab.Color = 11957550

# --- Active cell / selection moved to D2 ---
$ws.Range("D2").Select()
